# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") values for rows 2-5 are regenerated; strike counts collapse to 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
